$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Variable Type column (H2:H4) from "Continuous" to "Discrete"
$ws.Range("H2:H4").Value = "Discrete"

# Row 2 (POX/C)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 130
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 0.06

# Row 3 (C/A)
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.13
$ws.Range("E3").Value = 0.07000000000000001
$ws.Range("F3").Value = 0.00005999999999999999

# Row 4 (POX/M)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.0013
$ws.Range("E4").Value = 0.0007
$ws.Range("F4").Value = 0.0000006
